$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 731, shifting existing rows 731:806 down to 732:807
$ws.Rows.Item(731).Insert()

# Populate the newly inserted row 731 with the new data record
$ws.Cells.Item(731, 1).Value = 4
$ws.Cells.Item(731, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(731, 3).Value = "Los Lagos"
$ws.Cells.Item(731, 4).Value = 45166
$ws.Cells.Item(731, 5).Value = 10
$ws.Cells.Item(731, 6).Value = 100112006
$ws.Cells.Item(731, 7).Value = "Repollo"
$ws.Cells.Item(731, 8).Value = "Crespo record"
$ws.Cells.Item(731, 9).Value = "Primera"
$ws.Cells.Item(731, 10).Value = 500
$ws.Cells.Item(731, 11).Value = 1500
$ws.Cells.Item(731, 12).Value = 1500
$ws.Cells.Item(731, 13).Value = 1500
$ws.Cells.Item(731, 14).Value = "$/unidad"
$ws.Cells.Item(731, 15).Value = "Región Metropolitana"
$ws.Cells.Item(731, 16).Value = 1500
$ws.Cells.Item(731, 17).Value = 1
$ws.Cells.Item(731, 18).Value = "Hortaliza"

$wb.Save()
